# Atualização automática de GRAVATAI.xlsx
#
# 1) Rename "Paineis DARQ"            -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Delete the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete this sheet" prompts.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
